$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 526
$ws1.Range("F4").Value = 1521
$ws1.Range("F8").Value = 154
$ws1.Range("F9").Value = 739
$ws1.Range("F12").Value = 328
$ws1.Range("F14").Value = 6400
$ws1.Range("F15").Value = 9
$ws1.Range("F20").Value = 15316
$ws1.Range("F25").Value = 11047
$ws1.Range("F27").Value = 4321
$ws1.Range("F28").Value = 237

# Sheet "全部类型" (All types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 526
$ws4.Range("F4").Value = 1521
$ws4.Range("F9").Value = 154
$ws4.Range("F10").Value = 739
$ws4.Range("F14").Value = 328
$ws4.Range("F17").Value = 6400
$ws4.Range("F18").Value = 9
$ws4.Range("F23").Value = 15316
$ws4.Range("F28").Value = 11047
$ws4.Range("F30").Value = 4321
$ws4.Range("F31").Value = 237
